$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Home Care Subsidy" program as row 28, following the same
# column layout used by the rest of the table:
#   A program_identifier   B program_name   C year   D category
#   E average_age_beneficiary   F short_description   I Papers   J Links
$ws.Range("A28").Value = "homeCareSubsidy"
$ws.Range("B28").Value = "Home Care Subsidy"
$ws.Range("C28").Value = 2013
$ws.Range("D28").Value = "Other"
$ws.Range("E28").Value = 31.36
$ws.Range("F28").Value = "The home care subsidy ""Betreuungsgeld"" was introduced in 2013 and was meant to compensate parents who did not make use of subsidised childcare."
$ws.Range("I28").Value = "Collischon et al. (2020)"
$ws.Range("J28").Value = "https://ideas.repec.org/p/iza/izadps/dp13271.html"

# Give the new row the same row height as its neighbours.
$ws.Rows.Item(28).RowHeight = 60

# Turn the URL in J28 into a real hyperlink, same as the rest of column J.
$ws.Hyperlinks.Add($ws.Range("J28"), "https://ideas.repec.org/p/iza/izadps/dp13271.html")

# Hyperlinks.Add() re-styles the cell with a fresh "Link" style; restore the
# same look-and-feel (font/underline + wrap) already used by the other
# hyperlink cells in column J by copying the format from J27.
$ws.Range("J27").Copy()
$ws.Range("J28").PasteSpecial(-4122)

# Reflect where the user ended up after entering the new row.
$ws.Range("J28").Select()
